$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 ("Save") should match styling of G1 (bold, centered, bordered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells H2:H5
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
